$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record for Zapallo ("Paine" variety, "1a (guarda)") needs to be
# inserted into the table at row 113, pushing the existing rows 113-136 down to
# 114-137 (dimension grows from R136 to R137).
$ws.Rows.Item(113).Insert()

$ws.Range("A113").Value = 11
$ws.Range("B113").Value = "Vega Monumental Concepción"
$ws.Range("C113").Value = "Bíobío"
$ws.Range("D113").Value = 44511
$ws.Range("E113").Value = 8
$ws.Range("F113").Value = 100112045
$ws.Range("G113").Value = "Zapallo"
$ws.Range("H113").Value = "Paine"
$ws.Range("I113").Value = "1a (guarda)"
$ws.Range("J113").Value = 200
$ws.Range("K113").Value = 150
$ws.Range("L113").Value = 150
$ws.Range("M113").Value = 150
$ws.Range("N113").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O113").Value = "Región del Maule"
$ws.Range("P113").Value = 150
$ws.Range("Q113").Value = 1
$ws.Range("R113").Value = "Hortaliza"
